$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.943.55"
$ws.Range("E2").Value = "  +1.32%  "
# Row 3
$ws.Range("D3").Value = "2.584.51"
$ws.Range("E3").Value = "  -0.37%  "
# Row 4
$ws.Range("E4").Value = "  +0.26%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.77"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.17%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.87"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.55%  "
# Row 7
$ws.Range("E7").Value = "  +0.08%  "
# Row 8
$ws.Range("E8").Value = "  -0.55%  "
# Row 9
$ws.Range("D9").Value = "2.595.21"
$ws.Range("E9").Value = "  -0.77%  "
# Row 10
$ws.Range("E10").Value = "  -2.27%  "
# Row 11
$ws.Range("E11").Value = "  -0.57%  "
# Row 12
$ws.Range("E12").Value = "  -0.06%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.135"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +2.97%  "
# Row 14
$ws.Range("D14").Value = "3.045.61"
$ws.Range("E14").Value = "  -0.31%  "
# Row 15
$ws.Range("D15").Value = "58.896.97"
$ws.Range("E15").Value = "  +1.26%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.51"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +0.21%  "
# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.30%  "
# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.579.23"
$ws.Range("E18").Value = "  -1.11%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.90"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.88%  "
# Row 20
$ws.Range("E20").Value = "  -1.10%  "
# Row 21
$ws.Range("E21").Value = "  -2.36%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.82%  "
# Row 23
$ws.Range("E23").Value = "  +0.06%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.12"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.33%  "
# Row 25
$ws.Range("E25").Value = "  +1.18%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.403"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.05%  "
# Row 27
$ws.Range("E27").Value = "  +0.23%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.01"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.22%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.04%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0720"
$ws.Range("E30").Value = "  -4.07%  "
# Row 31
$ws.Range("E31").Value = "  -5.58%  "
# Row 32
$ws.Range("E32").Value = "  -0.54%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.66"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.99%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.50"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.16%  "
# Row 35
$ws.Range("E35").Value = "  -1.79%  "
# Row 36
$ws.Range("E36").Value = "  -2.92%  "
# Row 37
$ws.Range("E37").Value = "  +1.44%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.47"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.22%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.823"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -2.67%  "
# Row 40
$ws.Range("E40").Value = "  -6.63%  "
# Row 41
$ws.Range("E41").Value = "  -1.07%  "
# Row 42
$ws.Range("E42").Value = "  +0.18%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "271.83"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.06%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.77"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +1.15%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.592"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.10%  "
# Row 47
$ws.Range("E47").Value = "  -1.71%  "
# Row 48
$ws.Range("E48").Value = "  -2.98%  "
# Row 49
$ws.Range("D49").Value = "1.966.95"
$ws.Range("E49").Value = "  -1.17%  "
# Row 50
$ws.Range("E50").Value = "  -0.62%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.48"
$ws.Range("D51").NumberFormat = "General"
